# Adds a "Team" column (AF) to the Thomas Tuchel sheet, indicating which
# club Tuchel was coaching for in each fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell AF1, styled like the other header cells (e.g. copy from AE1)
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("AF1").Value = "Team"

# Data rows: which team Tuchel coached for, by row range
$ws.Range("AF2:AF171").Value = "Mainz"
$ws.Range("AF172:AF239").Value = "Dortmund"
$ws.Range("AF240:AF321").Value = "Paris SG"
$ws.Range("AF322:AF340").Value = "Chelsea"
